$wb = $excel.ActiveWorkbook

# --- Sheet "Linear" ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.00001559499027795207
$wsLinear.Range("B3").Value = -0.001973831047874729
$wsLinear.Range("B4").Value = 0.00172479833526934
$wsLinear.Range("B5").Value = "[1.0, 0.1525143732735503, -0.0071801845636677245, -0.06386388064652439, -0.08437758483280604, -0.026666015846221563, 0.15616924329792167, 0.3079629521290953, 0.12682109424832988, -0.03567518520866726, -0.10403850491525923, -0.07775068958576623, -0.05713840945199983, 0.1538747973035678, 0.3073728262932102, 0.13745764394194362, -0.03586596423707714, -0.10765067293781076, -0.09413323618878881, -0.04939421997027026]"

# --- Sheet "NonLinear" ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.00164860382354744
$wsNonLinear.Range("B5").Value = 0.001476759589827472
$wsNonLinear.Range("B6").Value = 0.001768294473691201
$wsNonLinear.Range("B7").Value = 0.0002090679752573348
$wsNonLinear.Range("B8").Value = -0.003020441998840973
$wsNonLinear.Range("B9").Value = 0.001683786525185128
$wsNonLinear.Range("B10").Value = "[1.0, 0.1523674428269201, -0.005535429335141594, -0.061912761488983566, -0.08328815727985144, -0.02595889001695908, 0.15494943008275708, 0.30557787031226324, 0.12560443642927338, -0.03441670491204005, -0.10249552638206944, -0.07632453255560184, -0.05645132754036239, 0.1524549843567909, 0.3051260436007956, 0.1365502049137924, -0.03515224251559686, -0.10548666175929189, -0.09298151434264713, -0.04906930652511911]"

Write-Host "Updated Linear and NonLinear calibration values"
